# Update countries & provincias Spain
# Applies the 29-Apr-2020 13:22 data refresh: several countries swapped rank
# order (so their row keeps the old country but gets new numbers, while the
# adjacent row receives the country that overtook/fell behind it) and the
# "Datos actualizados" timestamp moves from 12:52 to 13:22.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- timestamp banner (A1) -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 29 de Abril de 2020 a las 13:22"

# --- helper: write country name + the 7 numeric columns for a row ----------
function Set-Row($row, $country, $casosTotales, $nuevosCasos, $casosActivos, $recuperados, $casosCriticos, $muertesHoy, $muertes) {
    $ws.Range("A$row").Value = $country
    $ws.Range("B$row").Value = $casosTotales
    $ws.Range("C$row").Value = $nuevosCasos
    $ws.Range("D$row").Value = $casosActivos
    $ws.Range("E$row").Value = $recuperados
    $ws.Range("F$row").Value = $casosCriticos
    $ws.Range("G$row").Value = $muertesHoy
    $ws.Range("H$row").Value = $muertes
}

# Row 29 - Austria (no rank change, figures refreshed)
Set-Row 29 "Austria" 15402 45 12779 2043 131 11 580

# Rows 33/34 - Bielorrusia overtakes Polonia
Set-Row 33 "Bielorrusia" 13181 973 2072 11025 92 5 84
Set-Row 34 "Polonia" 12415 197 3025 8784 160 10 606

# Row 35 - Rumania (no rank change, figures refreshed)
Set-Row 35 "Rumania" 11978 362 3569 7728 247 18 681

# Rows 75/76/77 - Bosnia y Herzegovina jumps ahead of Ghana and Estonia
Set-Row 75 "Bosnia y Herzegovina" 1677 92 710 902 4 2 65
Set-Row 76 "Ghana" 1671 0 188 1467 4 0 16
Set-Row 77 "Estonia" 1666 6 236 1380 10 0 50

# Rows 106/107 - San Marino overtakes Guatemala
Set-Row 106 "San Marino" 563 10 69 453 6 0 41
Set-Row 107 "Guatemala" 557 27 62 479 5 1 16

# Rows 111/112 - Malta overtakes Mayotte
Set-Row 111 "Malta" 463 5 339 120 1 0 4
Set-Row 112 "Mayotte" 460 0 235 221 4 0 4

# Row 170 - Macao (no rank change, figures refreshed)
Set-Row 170 "Macao" 45 0 34 11 1 0 0
